# Apply the "build : added some changes" edit:
#   1. Remove the "id" column (column H) entirely.
#   2. Fix the typo'd birth date for row 2 (soheil): 09/0/2018' -> 09/06/2018'.
#   3. Append a new data row (hojat / sarvar) with its own hyperlink.
#   4. Update the sheet selection to D7 (matches the saved workbook view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the trailing "id" column (H) - shifts nothing else, just removes it.
$ws.Columns("H").Delete()

# 2) Correct the birthDayDate value for soheil (row 2), keeping its existing
#    date-ish display format.
$ws.Range("D2").Value = "09/06/2018'"
$ws.Range("D2").NumberFormat = "mm-dd-yy"

# 3) Add the new row (row 4) with hojat / sarvar's data.
$ws.Range("A4").Value = "hojat"
$ws.Range("B4").Value = "sarvar"
$ws.Range("C4").Value = 1203659875
$ws.Range("D4").Value = "02/07/2015'"
$ws.Range("E4").Value = 6543653
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = "hojat@gmail.com"
$ws.Range("G4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:hojat@gmail.com")

# 4) Match the saved selection/active cell from the source workbook.
$ws.Range("D7").Select()
